$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column K: intervention_type
# Copy the header formatting (bold, border, centered) from J1, then set the text
$ws.Range("J1").Copy($ws.Range("K1"))
$ws.Range("K1").Value = "intervention_type"

$ws.Range("K2").Value = "PROCEDURE"
$ws.Range("K3").Value = "OTHER"
$ws.Range("K4").Value = "DEVICE"
$ws.Range("K5").Value = "OTHER"
$ws.Range("K6").Value = "BIOLOGICAL"

# Row 7 keeps an empty (blank) cell in column K
$ws.Range("K7").Style = "Normal"
